$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New customer row 18: phone 79174449, no birthday on file, 0 points.
# Phone numbers in column A are stored as text in this sheet's new row
# (leading apostrophe forces text storage so the value round-trips as
# "79174449" rather than being coerced to a number).
$ws.Range("A18").Value = "'79174449"
$ws.Range("B18").Value = "'"
$ws.Range("C18").Value = 0
